$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the two blank placeholder rows (15 & 16) that sat above the old
#    "2nd Week" header; this slides that header up to row 15 (matching the
#    target layout) while reusing its existing style untouched.
# ---------------------------------------------------------------------------
$ws.Rows("15:16").Delete()

# ---------------------------------------------------------------------------
# 2. Open up 4 fresh rows (16-19) right below the header for the new
#    "2nd Week" task block.
# ---------------------------------------------------------------------------
$ws.Rows("16:19").Insert()

# ---------------------------------------------------------------------------
# 3. Row 16: category/header strip "ایلام ـ انجام بند 1 نامه ـ بخش دوم"
#    merged across A16:C16, same light-blue fill + right alignment used by
#    the other category rows (A9, A12).
# ---------------------------------------------------------------------------
$catFill = 15652797   # theme Accent1, 60% lighter (same as existing category rows)

$a16 = $ws.Cells.Item(16, 1)
$b16 = $ws.Cells.Item(16, 2)
$c16 = $ws.Cells.Item(16, 3)

foreach ($cell in @($a16, $b16, $c16)) {
    $cell.Interior.Color = $catFill
    $cell.HorizontalAlignment = -4152   # xlRight
    $cell.Borders(8).LineStyle = 1      # xlEdgeTop
    $cell.Borders(9).LineStyle = 1      # xlEdgeBottom
}
$a16.Borders(7).LineStyle = 1           # xlEdgeLeft
$c16.Borders(10).LineStyle = 1          # xlEdgeRight

$a16.Value = "ایلام ـ انجام بند 1 نامه ـ بخش دوم"
$ws.Range("A16:C16").Merge()

# ---------------------------------------------------------------------------
# 4. Rows 17-19: new task block details.
#    Column A (merged A17:A19) is a dark-blue / yellow-bold "app" banner;
#    columns B & C hold the individual file/report names using the normal
#    bordered style already used throughout the sheet.
# ---------------------------------------------------------------------------
$darkBlue = 11957550   # theme Accent1, 25% darker

$a17 = $ws.Cells.Item(17, 1)
$a18 = $ws.Cells.Item(18, 1)
$a19 = $ws.Cells.Item(19, 1)

foreach ($cell in @($a17, $a18, $a19)) {
    $cell.Interior.Color = $darkBlue
    $cell.Font.Color = 65535            # yellow
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4108     # xlCenter
    $cell.Borders(7).LineStyle = 1      # xlEdgeLeft
    $cell.Borders(10).LineStyle = 1     # xlEdgeRight
}
$a17.Borders(8).LineStyle = 1           # xlEdgeTop
$a19.Borders(9).LineStyle = 1           # xlEdgeBottom

$ws.Cells.Item(17, 2).Value = "frmMPFeederPeakDayNight.vb"
$ws.Cells.Item(19, 2).Value = "Report_8_24.mrt"
$ws.Cells.Item(17, 3).Value = "Report_8_24"
$a17.Value = "Havades_App"
$ws.Range("A17:A19").Merge()
$ws.Cells.Item(18, 2).Value = "frmMain.vb"

# ---------------------------------------------------------------------------
# 5. Two extra blank rows (25 & 26) appended at the bottom, matching the
#    plain bordered look of the rest of the table body.
# ---------------------------------------------------------------------------
for ($r = 25; $r -le 26; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Borders(7).LineStyle = 1
        $cell.Borders(8).LineStyle = 1
        $cell.Borders(9).LineStyle = 1
        $cell.Borders(10).LineStyle = 1
    }
}

# ---------------------------------------------------------------------------
# 6. Update the view state: scrolled down a bit with B23 selected, as in the
#    saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("A4").Select()
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B23").Select()
